$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.00469825789230973
$ws.Range("J2").Value = 0.00469825789230973
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 15.70622098470778
$ws.Range("R2").Value = 141.35598886237
$ws.Range("S2").Value = 0.001402045784589196
$ws.Range("T2").Value = 0.001402045784589196
$ws.Range("I3").Value = 0.00469825789230973
$ws.Range("J3").Value = 0.00469825789230973
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("S3").Value = 0.001359481531826468
$ws.Range("T3").Value = 0.001359481531826468
$ws.Range("I4").Value = 0.00469825789230973
$ws.Range("J4").Value = 0.00469825789230973
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 15.50849885201667
$ws.Range("R4").Value = 139.57648966815
$ws.Range("S4").Value = 0.001384395741149118
$ws.Range("T4").Value = 0.001384395741149118
$ws.Range("I5").Value = 0.00469825789230973
$ws.Range("J5").Value = 0.00469825789230973
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 6.187453410872778
$ws.Range("R5").Value = 55.687080697855
$ws.Range("S5").Value = 0.000552334834744949
$ws.Range("T5").Value = 0.000552334834744949
$ws.Range("I6").Value = 0.7185612021237531
$ws.Range("J6").Value = 0.7185612021237531
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 2402.141664055107
$ws.Range("R6").Value = 21619.27497649596
$ws.Range("S6").Value = 0.2144317590688223
$ws.Range("T6").Value = 0.2144317590688223
$ws.Range("I7").Value = 0.7185612021237531
$ws.Range("J7").Value = 0.7185612021237531
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.2079218949162505
$ws.Range("T7").Value = 0.2079218949162505
$ws.Range("I8").Value = 0.7185612021237531
$ws.Range("J8").Value = 0.7185612021237531
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 2371.90163538713
$ws.Range("R8").Value = 21347.11471848417
$ws.Range("S8").Value = 0.21173232521003
$ws.Range("T8").Value = 0.21173232521003
$ws.Range("I9").Value = 0.7185612021237531
$ws.Range("J9").Value = 0.7185612021237531
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 946.3218203238545
$ws.Range("R9").Value = 8516.896382914689
$ws.Range("S9").Value = 0.08447522292865027
$ws.Range("T9").Value = 0.08447522292865027
$ws.Range("G10").Value = 5.503190333333333
$ws.Range("H10").Value = 16.509571
$ws.Range("I10").Value = 0.2767405399839373
$ws.Range("J10").Value = 0.2767405399839373
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 925.1403767191358
$ws.Range("R10").Value = 8326.263390472222
$ws.Range("S10").Value = 0.08258442094984046
$ws.Range("T10").Value = 0.08258442094984048
$ws.Range("G11").Value = 5.503190333333333
$ws.Range("H11").Value = 16.509571
$ws.Range("I11").Value = 0.2767405399839373
$ws.Range("J11").Value = 0.2767405399839373
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 897.0543403938552
$ws.Range("R11").Value = 8073.489063544697
$ws.Range("S11").Value = 0.08007726732746254
$ws.Range("T11").Value = 0.08007726732746254
$ws.Range("G12").Value = 5.503190333333333
$ws.Range("H12").Value = 16.509571
$ws.Range("I12").Value = 0.2767405399839373
$ws.Range("J12").Value = 0.2767405399839373
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 913.4939896918767
$ws.Range("R12").Value = 8221.445907226891
$ws.Range("S12").Value = 0.08154478398986384
$ws.Range("T12").Value = 0.08154478398986384
$ws.Range("G13").Value = 5.503190333333333
$ws.Range("H13").Value = 16.509571
$ws.Range("I13").Value = 0.2767405399839373
$ws.Range("J13").Value = 0.2767405399839373
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 364.4583241914348
$ws.Range("R13").Value = 3280.124917722913
$ws.Range("S13").Value = 0.03253406771677044
$ws.Range("T13").Value = 0.03253406771677044
